$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.75
$ws.Range("I3").Value = 3.1
$ws.Range("J3").Value = 3.75
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 5
$ws.Range("X3").Value = 11
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 13
$ws.Range("AR3").Value = 151
$ws.Range("AX3").Value = 21

# Row 4
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = 2.2
$ws.Range("K4").Value = 2.2
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("X4").Value = 6.5
$ws.Range("Y4").Value = 8.5
$ws.Range("AC4").Value = 8.5
$ws.Range("AE4").Value = 21
$ws.Range("AJ4").Value = 21
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 8
$ws.Range("AQ4").Value = 26
$ws.Range("AU4").Value = 9.5
$ws.Range("AW4").Value = 7.5

# Row 6
$ws.Range("I6").Value = 4.75
$ws.Range("AG6").Value = 900
$ws.Range("AN6").Value = 3.6
$ws.Range("AO6").Value = 10

# Row 7
$ws.Range("I7").Value = 6.6
$ws.Range("P7").Value = 4.7
$ws.Range("AM7").Value = 45
$ws.Range("AU7").Value = 7.3
$ws.Range("AV7").Value = 50

# Row 8
$ws.Range("G8").Value = 2.32
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 2.9
$ws.Range("J8").Value = 2.92
$ws.Range("K8").Value = 2.07
$ws.Range("N8").Value = 7.9
$ws.Range("O8").Value = 1.24
$ws.Range("P8").Value = 3.65
$ws.Range("Q8").Value = 1.75
$ws.Range("R8").Value = 2.02
$ws.Range("T8").Value = 2.7
$ws.Range("W8").Value = 9.5
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 7.9
$ws.Range("AD8").Value = 6.5
$ws.Range("AK8").Value = 37
$ws.Range("AO8").Value = 12.5
$ws.Range("AP8").Value = 19.5
$ws.Range("AR8").Value = 80
$ws.Range("AT8").Value = 2.7
$ws.Range("AU8").Value = 6.8
$ws.Range("AV8").Value = 60

# Row 9
$ws.Range("G9").Value = 2.07
$ws.Range("I9").Value = 3.25
$ws.Range("L9").Value = 3.65
$ws.Range("P9").Value = 3.7
$ws.Range("T9").Value = 2.9
$ws.Range("W9").Value = 8.75
$ws.Range("X9").Value = 11
$ws.Range("Z9").Value = 19.5
$ws.Range("AB9").Value = 23
$ws.Range("AH9").Value = 11.75
$ws.Range("AI9").Value = 18.5
$ws.Range("AK9").Value = 45
$ws.Range("AM9").Value = 29
$ws.Range("AN9").Value = 4.1
$ws.Range("AO9").Value = 10.5
$ws.Range("AT9").Value = 2.9
$ws.Range("AU9").Value = 6.7
$ws.Range("AW9").Value = 5.3
$ws.Range("AY9").Value = 22
$ws.Range("BA9").Value = 100

# Row 10
$ws.Range("G10").Value = 2.72
$ws.Range("H10").Value = 3.2
$ws.Range("J10").Value = 3.3
$ws.Range("L10").Value = 3.15
$ws.Range("Q10").Value = 1.88
$ws.Range("W10").Value = 9.25
$ws.Range("Y10").Value = 9.75
$ws.Range("AA10").Value = 22
$ws.Range("AB10").Value = 28
$ws.Range("AD10").Value = 6.2
$ws.Range("AH10").Value = 8.5
$ws.Range("AI10").Value = 12.5
$ws.Range("AK10").Value = 27
$ws.Range("AL10").Value = 20
$ws.Range("AM10").Value = 28
$ws.Range("AN10").Value = 4.7
$ws.Range("AO10").Value = 14.5
$ws.Range("AP10").Value = 21
$ws.Range("AQ10").Value = 65
$ws.Range("AR10").Value = 90
$ws.Range("AU10").Value = 6.8
$ws.Range("AV10").Value = 60
$ws.Range("AX10").Value = 14
$ws.Range("AY10").Value = 21
$ws.Range("AZ10").Value = 60
$ws.Range("BA10").Value = 100
